# Auto-generated PowerShell COM-interop script
# Adds column G (2025/11/13 data) to StockData sheet, mirroring the
# existing layout/formatting used by column F (and E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Apply formatting to column G, matching the pattern used by column F ----

# Numeric data rows (style matching F-column numeric cells: 0.00 format, centered, 12pt)
$rng = $ws.Range("G3:G4")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G6:G7")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G9:G10")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G12:G13")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G15:G16")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G18:G19")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G21:G22")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G24:G25")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G27:G28")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G30:G31")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G33:G34")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G36:G37")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G39:G40")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G42:G43")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G45:G46")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G48:G49")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G51:G52")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G54:G55")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G57:G58")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G60:G61")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G63:G64")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G66:G67")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G69:G70")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G72:G73")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G75:G76")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G78:G79")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G81:G82")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G84:G85")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G87:G88")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G90:G91")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G93:G94")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G96:G97")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G99:G100")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G102:G103")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G105:G106")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G108:G109")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G111:G112")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G114:G115")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "0.00"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Blank separator rows (style matching F-column blank cells: General format, centered, 12pt)
$rng = $ws.Range("G5")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G8")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G11")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G14")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G17")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G20")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G23")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G26")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G29")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G32")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G35")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G38")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G41")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G44")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G47")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G50")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G53")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G56")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G59")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G62")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G65")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G68")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G71")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G74")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G77")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G80")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G83")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G86")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G89")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G92")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G95")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G98")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G101")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G104")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G107")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G110")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$rng = $ws.Range("G113")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Header date row (style matching F1: General text, centered, 12pt, not bold)
$rng = $ws.Range("G1")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.NumberFormat = "@"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# Label row (style matching F2: bold text, centered, 12pt)
$rng = $ws.Range("G2")
$rng.Font.Name = "宋体"
$rng.Font.Size = 12
$rng.Font.Bold = $true
$rng.NumberFormat = "General"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# ---- 2) Write the values for each row of the new column G ----

$ws.Range("G1").Value = "2025/11/13"
$ws.Range("G2").Value = "上证"
$ws.Range("G3").Value = 63.21
$ws.Range("G4").Value = 4017.94
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 5648.68
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 4690.71
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 7363.44
$ws.Range("G15").Value = 29.44
$ws.Range("G16").Value = 2755.62
$ws.Range("G18").Value = 96.67
$ws.Range("G19").Value = 6850.92
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 84525.89
$ws.Range("G24").Value = 84.14
$ws.Range("G25").Value = 19909.14
$ws.Range("G27").Value = 71.64
$ws.Range("G28").Value = 39894.54
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = 5774
$ws.Range("G33").Value = 12.61
$ws.Range("G34").Value = 32685.54
$ws.Range("G36").Value = 28.23
$ws.Range("G37").Value = 3381.23
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 3205.76
$ws.Range("G42").Value = 18.21
$ws.Range("G43").Value = 7305.09
$ws.Range("G45").Value = 0
$ws.Range("G46").Value = 8880.959999999999
$ws.Range("G48").Value = 0
$ws.Range("G49").Value = 13055.99
$ws.Range("G51").Value = 26.04
$ws.Range("G52").Value = 12603.71
$ws.Range("G54").Value = 0
$ws.Range("G55").Value = 10004.5
$ws.Range("G57").Value = 0
$ws.Range("G58").Value = 16387.97
$ws.Range("G60").Value = 0
$ws.Range("G61").Value = 17526.85
$ws.Range("G63").Value = 21.11
$ws.Range("G64").Value = 10214.48
$ws.Range("G66").Value = 19.06
$ws.Range("G67").Value = 10219.1
$ws.Range("G69").Value = 0
$ws.Range("G70").Value = 3100.16
$ws.Range("G72").Value = 0
$ws.Range("G73").Value = 5886.87
$ws.Range("G75").Value = 0
$ws.Range("G76").Value = 9468.610000000001
$ws.Range("G78").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("G81").Value = 56.04
$ws.Range("G82").Value = 3076.44
$ws.Range("G84").Value = 0
$ws.Range("G85").Value = 3158.76
$ws.Range("G87").Value = 51.86
$ws.Range("G88").Value = 4117.83
$ws.Range("G90").Value = 0
$ws.Range("G91").Value = 2075.4
$ws.Range("G93").Value = 28.25
$ws.Range("G94").Value = 14198.75
$ws.Range("G96").Value = 85.78
$ws.Range("G97").Value = 9121.799999999999
$ws.Range("G99").Value = 56.31
$ws.Range("G100").Value = 12067.47
$ws.Range("G102").Value = 6.14
$ws.Range("G103").Value = 2306.9
$ws.Range("G105").Value = 25.81
$ws.Range("G106").Value = 873.11
$ws.Range("G108").Value = 0
$ws.Range("G109").Value = 2925.23
$ws.Range("G111").Value = 0
$ws.Range("G112").Value = 4021.92
$ws.Range("G114").Value = 29.02
$ws.Range("G115").Value = 3410

# ---- 3) Update the sheet dimension to reflect the new column ----
Write-Host "Column G populated (rows 1-115)."